# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across all job
# sheets in the Brynhildr_Profits workbook.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets("ALC")
$ws.Range("H62").Value = 3561.8235
$ws.Range("I62").Value = 3561.8235
$ws.Range("K62").Value = 3561.8235
$ws.Range("M62").Value = -2937.8235
$ws.Range("H65").Value = 3561.8235
$ws.Range("I65").Value = 3561.8235
$ws.Range("K65").Value = 17809.1175
$ws.Range("M65").Value = -14689.1175
$ws.Range("H80").Value = 1494.125
$ws.Range("I80").Value = 878.2222
$ws.Range("J80").Value = 2286
$ws.Range("K80").Value = 2634.6666
$ws.Range("L80").Value = 6858
$ws.Range("M80").Value = -1636.6666
$ws.Range("N80").Value = -8854
$ws.Range("H83").Value = 1494.125
$ws.Range("I83").Value = 878.2222
$ws.Range("J83").Value = 2286
$ws.Range("K83").Value = 7903.999800000001
$ws.Range("L83").Value = 20574
$ws.Range("M83").Value = -2911.999800000001
$ws.Range("N83").Value = -30558
$ws.Range("H100").Value = 4531
$ws.Range("I100").Value = 1651.75
$ws.Range("J100").Value = 6066.6
$ws.Range("K100").Value = 1651.75
$ws.Range("L100").Value = 6066.6
$ws.Range("M100").Value = -1110.75
$ws.Range("N100").Value = -7148.6
$ws.Range("H113").Value = 3999
$ws.Range("I113").Value = 3666.6667
$ws.Range("K113").Value = 3666.6667
$ws.Range("M113").Value = -412.6667000000002

# ----- Sheet: ARM -----
$ws = $wb.Worksheets("ARM")
$ws.Range("H74").Value = 5001.4443
$ws.Range("I74").Value = 2685.6345
$ws.Range("K74").Value = 2685.6345
$ws.Range("M74").Value = -1811.6345
$ws.Range("H77").Value = 5001.4443
$ws.Range("I77").Value = 2685.6345
$ws.Range("K77").Value = 13428.1725
$ws.Range("M77").Value = -9060.172500000001
$ws.Range("H132").Value = 4766.127
$ws.Range("I132").Value = 2820.7021
$ws.Range("K132").Value = 8462.106299999999
$ws.Range("M132").Value = -5932.106299999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets("BSM")
$ws.Range("H82").Value = 30534.715
$ws.Range("I82").Value = 18185.75
$ws.Range("K82").Value = 18185.75
$ws.Range("M82").Value = -17802.75
$ws.Range("H85").Value = 30534.715
$ws.Range("I85").Value = 18185.75
$ws.Range("K85").Value = 18185.75
$ws.Range("M85").Value = -16859.75
$ws.Range("H94").Value = 4995.4546
$ws.Range("I94").Value = 4913.684
$ws.Range("J94").Value = 5513.3335
$ws.Range("K94").Value = 4913.684
$ws.Range("L94").Value = 5513.3335
$ws.Range("M94").Value = -4462.684
$ws.Range("N94").Value = -6415.3335

# ----- Sheet: CRP -----
$ws = $wb.Worksheets("CRP")
$ws.Range("H122").Value = 12969.583
$ws.Range("I122").Value = 2437.2222
$ws.Range("J122").Value = 44566.668
$ws.Range("K122").Value = 7311.6666
$ws.Range("L122").Value = 133700.004
$ws.Range("M122").Value = -4861.6666
$ws.Range("N122").Value = -138600.004
$ws.Range("H134").Value = 2144.96
$ws.Range("I134").Value = 1883.6471
$ws.Range("K134").Value = 5650.9413
$ws.Range("M134").Value = -3115.9413

# ----- Sheet: CUL -----
$ws = $wb.Worksheets("CUL")
$ws.Range("H3").Value = 10900
$ws.Range("I3").Value = 3700
$ws.Range("J3").Value = 14500
$ws.Range("K3").Value = 11100
$ws.Range("L3").Value = 43500
$ws.Range("M3").Value = -10988
$ws.Range("N3").Value = -43724
$ws.Range("H33").Value = 1443.6
$ws.Range("I33").Value = 112.333336
$ws.Range("J33").Value = 2532.818
$ws.Range("K33").Value = 674.000016
$ws.Range("L33").Value = 15196.908
$ws.Range("M33").Value = -391.000016
$ws.Range("N33").Value = -15762.908
$ws.Range("H96").Value = 13965
$ws.Range("I96").Value = 9825
$ws.Range("K96").Value = 29475
$ws.Range("M96").Value = -27416
$ws.Range("H113").Value = 20135.129
$ws.Range("I113").Value = 506.85715
$ws.Range("K113").Value = 1520.57145
$ws.Range("M113").Value = 649.4285500000001
$ws.Range("H137").Value = 9538.786
$ws.Range("J137").Value = 9885.615
$ws.Range("L137").Value = 29656.845
$ws.Range("N137").Value = -39856.845
$ws.Range("H138").Value = 24689
$ws.Range("I138").Value = 42310
$ws.Range("J138").Value = 13845.308
$ws.Range("K138").Value = 126930
$ws.Range("L138").Value = 41535.924
$ws.Range("M138").Value = -121790
$ws.Range("N138").Value = -51815.924

# ----- Sheet: GSM -----
$ws = $wb.Worksheets("GSM")
$ws.Range("H41").Value = 2280.4
$ws.Range("I41").Value = 2280.4
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2280.4
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1925.4
$ws.Range("N41").ClearContents()
$ws.Range("H113").Value = 1901.5333
$ws.Range("I113").Value = 1921.7084
$ws.Range("K113").Value = 1921.7084
$ws.Range("M113").Value = 248.2916
$ws.Range("H132").Value = 15686.963
$ws.Range("I132").Value = 20720.422
$ws.Range("J132").Value = 3732.5
$ws.Range("K132").Value = 62161.266
$ws.Range("L132").Value = 11197.5
$ws.Range("M132").Value = -59631.266
$ws.Range("N132").Value = -16257.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 3433.3125
$ws.Range("I7").Value = 3103.9
$ws.Range("J7").Value = 3982.3333
$ws.Range("K7").Value = 3103.9
$ws.Range("L7").Value = 3982.3333
$ws.Range("M7").Value = -2991.9
$ws.Range("N7").Value = -4206.3333
$ws.Range("H22").Value = 2524.5
$ws.Range("I22").Value = 2010
$ws.Range("K22").Value = 2010
$ws.Range("M22").Value = -1715
$ws.Range("H27").Value = 2524.5
$ws.Range("I27").Value = 2010
$ws.Range("K27").Value = 2010
$ws.Range("M27").Value = -1903
$ws.Range("H61").Value = 12516.728
$ws.Range("I61").Value = 12768.4
$ws.Range("K61").Value = 12768.4
$ws.Range("M61").Value = -12566.4
$ws.Range("H113").Value = 12516.728
$ws.Range("I113").Value = 12768.4
$ws.Range("K113").Value = 12768.4
$ws.Range("M113").Value = -10598.4
$ws.Range("H126").Value = 3433.3125
$ws.Range("I126").Value = 3103.9
$ws.Range("J126").Value = 3982.3333
$ws.Range("K126").Value = 9311.700000000001
$ws.Range("L126").Value = 11946.9999
$ws.Range("M126").Value = -6841.700000000001
$ws.Range("N126").Value = -16886.9999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets("WVR")
$ws.Range("H54").Value = 353348.34
$ws.Range("H81").Value = 22780674
$ws.Range("I81").Value = 3762.0833
$ws.Range("J81").Value = 50112970
$ws.Range("K81").Value = 7524.1666
$ws.Range("L81").Value = 100225940
$ws.Range("M81").Value = -6463.1666
$ws.Range("N81").Value = -100228062
$ws.Range("H84").Value = 22780674
$ws.Range("I84").Value = 3762.0833
$ws.Range("J84").Value = 50112970
$ws.Range("K84").Value = 37620.833
$ws.Range("L84").Value = 501129700
$ws.Range("M84").Value = -32316.833
$ws.Range("N84").Value = -501140308
$ws.Range("H132").Value = 3061.7222
$ws.Range("I132").Value = 2735.92
$ws.Range("K132").Value = 8207.76
$ws.Range("M132").Value = -5677.76
$ws.Range("H136").Value = 1935.9474
$ws.Range("I136").Value = 1660.2142
$ws.Range("K136").Value = 4980.642599999999
$ws.Range("M136").Value = -2430.642599999999
$ws.Range("H140").Value = 81250
$ws.Range("J140").Value = 81250
$ws.Range("L140").Value = 81250
$ws.Range("N140").Value = -91610
